$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.878.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.379.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.16%  "
$ws.Range("E7").Value = "  -3.92%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0841"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.28%  "
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.745.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.380.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.765"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "40.822.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0913"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("E24").Value = "  -3.68%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.92%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0733"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.23%  "
$ws.Range("E39").Value = "  -4.10%  "
$ws.Range("E40").Value = "  -7.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.958.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0272"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.51%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.604.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "93.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.10%  "
